$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert two new columns before column B, shifting old B:T -> D:V
$ws.Range("B:C").EntireColumn.Insert()

# --- Row 1 (span(m) =) new values ---
$ws.Range("B1").Value = 0.0126
$ws.Range("C1").Value = 0.025

# --- Row 2 (span(r/R) =) new formulas ---
$ws.Range("B2").Formula = "=B1/$A8"
$ws.Range("C2").Formula = "=C1/$A8"

# --- Row 3 (chord(m) =) new values ---
$ws.Range("B3").Value = 0.024
$ws.Range("C3").Value = 0.025

# --- Row 4 (Beta(deg) =) new values ---
$ws.Range("B4").Value = 15
$ws.Range("C4").Value = 28

# --- Row 5 (Foil Number =) new values ---
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1

# --- Row 6 (Delta r =) new formulas ---
$ws.Range("B6").Formula = "=B1-B8"
$ws.Range("C6").Formula = "=C1-B1"
